# "se actualiza la data para pagos"
# Update the payment subscription test data (numeroUsuario, placa, vigencia)
# across the sheets/tables that carry those columns. Leading apostrophes
# force text storage (matching the workbook's quotePrefix text styles),
# so the purely-numeric/date-looking values don't get reinterpreted as a
# number or a date serial by Excel.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("SucripcionDesdePagos")
$ws1.Range("A2").Value = "'72636759"
$ws1.Range("F2").Value = "'ZIJ-583"
$ws1.Range("G2").Value = "'12/11/2025"

$ws3 = $wb.Worksheets.Item("SucripcionDesdeAfiliacion")
$ws3.Range("A2").Value = "'72636759"
$ws3.Range("F2").Value = "'ZIJ-583"
$ws3.Range("G2").Value = "'12/11/2025"

$ws4 = $wb.Worksheets.Item("SucripcionDesdeOpcionPagar")
$ws4.Range("A2").Value = "'72636759"
$ws4.Range("C2").Value = "'ZIJ-583"
$ws4.Range("G2").Value = "'12/11/2025"
